# Thermal_Conductivity_Calculator.xlsx edit
#
# Summary of the change (per the commit's canonical-OOXML diff):
#   - On "Thermal Conductivity Values", 25 elements (technetium, promethium,
#     astatine, francium, neptunium, and most of the super-heavy/synthetic
#     elements from berkelium through oganesson) previously carried the
#     "unknown value" placeholder of 5000 in both the W/mK and BTU/hr-ft-F
#     columns. Those placeholders are replaced with real (researched)
#     thermal-conductivity figures.
#   - On "Calculator", the selected element (B42) changes from Neon to
#     Hydrogen, which ripples through the INDEX/MATCH formulas (C42, F42)
#     and the heat-flux table (C50, D50).
#   - Selection/active-cell bookkeeping is updated on both sheets.

$wb = $excel.ActiveWorkbook
$calc = $wb.Worksheets.Item("Calculator")
$vals = $wb.Worksheets.Item("Thermal Conductivity Values")

# --- Fill in the previously-unknown (5000 placeholder) thermal conductivity
#     values with real data -----------------------------------------------
$newValues = @{
    52  = @(50.6, 29.4)     # Tc  Technetium
    70  = @(13.9, 8.1)      # Pm  Promethium
    94  = @(0.0017, 0.001)  # At  Astatine
    96  = @(0.14, 0.0081)   # Fr  Francium
    102 = @(6, 3.49)        # Np  Neptunium
    106 = @(10, 5.78)       # Bk  Berkelium
    108 = @(10, 5.78)       # Es  Einsteinium
    109 = @(10, 5.78)       # Fm  Fermium
    110 = @(10, 5.78)       # Md  Mendelevium
    111 = @(10, 5.78)       # No  Nobelium
    112 = @(10, 5.78)       # Lr  Lawrencium
    114 = @(268, 156)       # Db  Dubnium
    115 = @(271, 158)       # Sg  Seaborgium
    116 = @(270, 157)       # Bh  Bohrium
    117 = @(277, 162)       # Hs  Hassium
    118 = @(278, 162)       # Mt  Meitnerium
    119 = @(281, 163)       # Ds  Darmstadtium
    120 = @(282, 165)       # Rg  Roentgenium
    121 = @(285, 167)       # Cn  Copernicium
    122 = @(284, 175)       # Nh  Nihonium
    123 = @(289, 170)       # Fl  Flerovium
    124 = @(288, 173)       # Mc  Moscovium
    125 = @(293, 174)       # Lv  Livermorium
    126 = @(294, 176)       # Ts  Tennessine
    127 = @(294, 176)       # Og  Oganesson
}

foreach ($row in $newValues.Keys) {
    $pair = $newValues[$row]
    $vals.Range("E$row").Value = $pair[0]
    $vals.Range("F$row").Value = $pair[1]
}

# E52 had a colour-scale conditional format that only existed to flag the
# "5000 == unknown" placeholder; now that it has a real value, drop the rule.
$vals.Range("E52").FormatConditions.Delete()

# --- Change the element selected in the calculator from Neon to Hydrogen -
$calc.Range("B42").Value = "Hydrogen"

# --- Restore/update the on-screen selection for each sheet ---------------
$vals.Activate()
$vals.Range("I16").Select()

$calc.Activate()
$calc.Range("B42").Select()
